$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 179; this shifts the existing rows 179:244 down to 180:245
$ws.Rows(179).Insert()

# Populate the newly inserted row 179 with its data
$ws.Cells.Item(179, 1).Value = 10
$ws.Cells.Item(179, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(179, 3).Value = "La Araucanía"
$ws.Cells.Item(179, 4).Value = 44559
$ws.Cells.Item(179, 5).Value = 9
$ws.Cells.Item(179, 6).Value = 100112044
$ws.Cells.Item(179, 7).Value = "Perejil"
$ws.Cells.Item(179, 8).Value = "Sin especificar"
$ws.Cells.Item(179, 9).Value = "Primera"
$ws.Cells.Item(179, 10).Value = 60
$ws.Cells.Item(179, 11).Value = 4500
$ws.Cells.Item(179, 12).Value = 5000
$ws.Cells.Item(179, 13).Value = 4708
$ws.Cells.Item(179, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(179, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(179, 16).Value = 1569
$ws.Cells.Item(179, 17).Value = 3
$ws.Cells.Item(179, 18).Value = "Hortaliza"
